# Refresh the cryptocurrency list: updated prices/volumes, and the
# HuobiToken / MultiversX rows (48-49) swapped places with new figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.211.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.309.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.51%  "
$ws.Range("E7").Value = "  +2.38%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0818"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.72%  "
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.666.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("E15").Value = "  +4.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.303.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("E17").Value = "  +3.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.115.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0921"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("E24").Value = "  +5.19%  "
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.08%  "
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0741"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("E37").Value = "  +3.20%  "
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.66%  "
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.986.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.20%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +18.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.533.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("E51").Value = "  +2.29%  "
